$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 39 ("Tolmie Peak" is row 38, "Van Trump..." was row 39).
# This shifts the existing rows 39-40 down to 40-41, keeping the table sorted
# alphabetically once we fill in "Twin Firs Loop".
$ws.Rows.Item(39).Insert()

# Grow the Excel Table (ListObject) so it covers the newly inserted row.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D41"))

# Populate the new row with the "Twin Firs Loop" hike data.
$ws.Range("A39").Value = "Twin Firs Loop"
$ws.Range("B39").Value = 0.4
$ws.Range("C39").Value = 180
$ws.Range("D39").Value = "easy"

# Match Excel's post-insert active selection (one cell below the last table row).
$ws.Range("D42").Select()
